$d = $word.ActiveDocument

# Find the "Abstract" title paragraph (style AbstractTitle), then remove the
# two "Block Text" styled paragraphs that immediately follow it (the
# instructional placeholder text), leaving the title paragraph directly
# followed by the real abstract content.

$count = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Abstract Title") {
        $j = $i + 1
        while ($j -le $count -and $d.Paragraphs.Item($j).Style.NameLocal -eq "Block Text") {
            $j = $j + 1
        }
        if ($j -gt ($i + 1)) {
            $startIndex = $i + 1
            $endIndex = $j - 1
        }
        break
    }
}

if ($startIndex -gt 0) {
    $start = $d.Paragraphs.Item($startIndex).Range.Start
    $end = $d.Paragraphs.Item($endIndex).Range.End
    $range = $d.Range($start, $end)
    $range.Delete()
}
